$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from January to November
$ws.Name = "November"

# New headers
$ws.Range("D1").Value = "People"
$ws.Range("E1").Value = "Nights"
$ws.Range("F1").Value = "Country"
$ws.Range("G1").Value = "Passport #"

# Column widths: D:E -> 10, F:G -> 15
$ws.Range("D1:E1").ColumnWidth = 10
$ws.Range("F1:G1").ColumnWidth = 15

# Row 2 data - 11/25/2024 as its Excel date serial number (avoids the
# engine auto-minting a throwaway date style when a true DateTime is
# assigned directly to .Value)
$ws.Range("A2").Value = 45621
$ws.Range("B2").Value = "1"
$ws.Range("C2").Value = "1"
$ws.Range("D2").Value = "1"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = "1"
$ws.Range("G2").Value = "1"

# Row 3 data
$ws.Range("A3").Value = 45621
$ws.Range("B3").Value = "STEF"
$ws.Range("C3").Value = "WG"
$ws.Range("D3").Value = "2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = "USA"
$ws.Range("G3").Value = "69"

# Date format for the Date column - apply to A2 then copy the format onto
# A3 so both cells share a single cell-format entry (avoids a duplicate
# style being minted for the second cell).
$ws.Range("A2").NumberFormat = "mm-dd-yy"
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = 45621
$excel.CutCopyMode = 0

# Page setup tweaks
$ws.PageSetup.FirstPageNumber = 1
$ws.PageSetup.Copies = 1
